# Update column F ("dSF") values for specific rows, per repull/push of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 3
    11 = 1
    12 = -3
    18 = 0
    20 = -3
    21 = 1
    23 = -1
    25 = -2
    31 = 1
    34 = 2
    36 = 6
    40 = -4
    44 = -1
    49 = -2
    55 = 5
    56 = 1
    63 = 1
    69 = 2
    73 = -3
    74 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
